# Auto-generated edit script applying updated profit/pricing values
# across ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR sheets.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 19
$ws.Range("H19").Value = 1893.2
$ws.Range("I19").Value = 2498.5
$ws.Range("J19").Value = 1489.6666
$ws.Range("K19").Value = 2498.5
$ws.Range("L19").Value = 1489.6666
$ws.Range("M19").Value = -2323.5
$ws.Range("N19").Value = -1839.6666

# Row 33
$ws.Range("H33").Value = 148.13333
$ws.Range("I33").Value = 148.13333
$ws.Range("K33").Value = 148.13333
$ws.Range("M33").Value = 80.86667

# Row 87
$ws.Range("H87").Value = 30000
$ws.Range("J87").Value = 30000
$ws.Range("L87").Value = 30000
$ws.Range("N87").Value = -32496

# Row 90
$ws.Range("H90").Value = 30000
$ws.Range("J90").Value = 30000
$ws.Range("L90").Value = 90000
$ws.Range("N90").Value = -102480

# Row 106
$ws.Range("H106").Value = 5500
$ws.Range("I106").Value = 5500
$ws.Range("K106").Value = 5500
$ws.Range("M106").Value = -4869

# Row 135
$ws.Range("H135").Value = 2040.3334
$ws.Range("I135").Value = 848.4
$ws.Range("K135").Value = 7635.599999999999
$ws.Range("M135").Value = -5100.599999999999

$ws = $wb.Worksheets.Item("ARM")
# Row 2
$ws.Range("H2").Value = 1923.1111
$ws.Range("I2").Value = 926.4286
$ws.Range("J2").Value = 5411.5
$ws.Range("K2").Value = 926.4286
$ws.Range("L2").Value = 5411.5
$ws.Range("M2").Value = -813.4286
$ws.Range("N2").Value = -5637.5

# Row 32
$ws.Range("H32").Value = 1771.1333
$ws.Range("I32").Value = 1388.4186
$ws.Range("K32").Value = 1388.4186
$ws.Range("M32").Value = -1101.4186

# Row 45
$ws.Range("H45").Value = 1502
$ws.Range("I45").Value = 1502
$ws.Range("K45").Value = 1502
$ws.Range("M45").Value = -1125

# Row 57
$ws.Range("H57").Value = 7875
$ws.Range("I57").Value = 7875
$ws.Range("K57").Value = 7875
$ws.Range("M57").Value = -7391

# Row 74
$ws.Range("H74").Value = 2852.2727
$ws.Range("I74").Value = 2852.2727
$ws.Range("K74").Value = 2852.2727
$ws.Range("M74").Value = -1978.2727

# Row 77
$ws.Range("H77").Value = 2852.2727
$ws.Range("I77").Value = 2852.2727
$ws.Range("K77").Value = 14261.3635
$ws.Range("M77").Value = -9893.363499999999

# Row 102
$ws.Range("H102").Value = 1987.5
$ws.Range("I102").Value = 1987.5
$ws.Range("K102").Value = 1987.5
$ws.Range("M102").Value = -365.5

# Row 110
$ws.Range("H110").Value = 1714.6923
$ws.Range("I110").Value = 1607.5834
$ws.Range("K110").Value = 1607.5834
$ws.Range("M110").Value = 437.4166

# Row 116
$ws.Range("H116").Value = 1923.1111
$ws.Range("I116").Value = 926.4286
$ws.Range("J116").Value = 5411.5
$ws.Range("K116").Value = 926.4286
$ws.Range("L116").Value = 5411.5
$ws.Range("M116").Value = 1367.5714
$ws.Range("N116").Value = -9999.5

# Row 122
$ws.Range("H122").Value = 2605.8125
$ws.Range("I122").Value = 2571
$ws.Range("K122").Value = 7713
$ws.Range("M122").Value = -5263

$ws = $wb.Worksheets.Item("BSM")
# Row 3
$ws.Range("H3").Value = 1923.1111
$ws.Range("I3").Value = 926.4286
$ws.Range("J3").Value = 5411.5
$ws.Range("K3").Value = 926.4286
$ws.Range("L3").Value = 5411.5
$ws.Range("M3").Value = -812.4286
$ws.Range("N3").Value = -5639.5

# Row 94
$ws.Range("H94").Value = 0
$ws.Range("I94").Value = 0
$ws.Range("K94").Value = 0
$ws.Range("M94").ClearContents()

# Row 99
$ws.Range("H99").Value = 4427.4287
$ws.Range("I99").Value = 4000
$ws.Range("J99").Value = 4997.3335
$ws.Range("K99").Value = 4000
$ws.Range("L99").Value = 4997.3335
$ws.Range("M99").Value = -2502
$ws.Range("N99").Value = -7993.3335

# Row 107
$ws.Range("H107").Value = 417.75
$ws.Range("I107").Value = 424.5
$ws.Range("K107").Value = 424.5
$ws.Range("M107").Value = 1495.5

# Row 134
$ws.Range("H134").Value = 5451.4136
$ws.Range("I134").Value = 7337.9287
$ws.Range("J134").Value = 3690.6667
$ws.Range("K134").Value = 22013.7861
$ws.Range("L134").Value = 11072.0001
$ws.Range("M134").Value = -19478.7861
$ws.Range("N134").Value = -16142.0001

$ws = $wb.Worksheets.Item("CRP")
# Row 16
$ws.Range("H16").Value = 1524
$ws.Range("I16").Value = 1524
$ws.Range("K16").Value = 1524
$ws.Range("M16").Value = -1237

# Row 31
$ws.Range("H31").Value = 1498
$ws.Range("I31").Value = 1308.8889
$ws.Range("J31").Value = 3200
$ws.Range("K31").Value = 1308.8889
$ws.Range("L31").Value = 3200
$ws.Range("M31").Value = -1013.8889
$ws.Range("N31").Value = -3790

# Row 34
$ws.Range("H34").Value = 1498
$ws.Range("I34").Value = 1308.8889
$ws.Range("J34").Value = 3200
$ws.Range("K34").Value = 1308.8889
$ws.Range("L34").Value = 3200
$ws.Range("M34").Value = -1106.8889
$ws.Range("N34").Value = -3604

# Row 92
$ws.Range("H92").Value = 55000
$ws.Range("J92").Value = 55000
$ws.Range("L92").Value = 55000
$ws.Range("N92").Value = -59992

# Row 113
$ws.Range("H113").Value = 1524
$ws.Range("I113").Value = 1524
$ws.Range("K113").Value = 1524
$ws.Range("M113").Value = 646

$ws = $wb.Worksheets.Item("CUL")
# Row 2
$ws.Range("H2").Value = 900
$ws.Range("I2").Value = 800
$ws.Range("J2").Value = 1000
$ws.Range("K2").Value = 4800
$ws.Range("L2").Value = 6000
$ws.Range("M2").Value = -4687
$ws.Range("N2").Value = -6226

# Row 39
$ws.Range("H39").Value = 2166.6667
$ws.Range("J39").Value = 2166.6667
$ws.Range("L39").Value = 6500.000100000001
$ws.Range("N39").Value = -7088.000100000001

# Row 40
$ws.Range("H40").Value = 101.21739
$ws.Range("I40").Value = 83
$ws.Range("J40").Value = 142.85715
$ws.Range("K40").Value = 332
$ws.Range("L40").Value = 571.4286
$ws.Range("M40").Value = -263
$ws.Range("N40").Value = -709.4286

# Row 68
$ws.Range("H68").Value = 1299.8
$ws.Range("J68").Value = 1750
$ws.Range("L68").Value = 5250
$ws.Range("N68").Value = -6872

# Row 71
$ws.Range("H71").Value = 1299.8
$ws.Range("J71").Value = 1750
$ws.Range("L71").Value = 15750
$ws.Range("N71").Value = -23862

# Row 107
$ws.Range("H107").Value = 999
$ws.Range("J107").Value = 999
$ws.Range("L107").Value = 2997
$ws.Range("N107").Value = -6837

# Row 122
$ws.Range("H122").Value = 619.2857
$ws.Range("I122").Value = 493.25
$ws.Range("J122").Value = 787.3333
$ws.Range("K122").Value = 4439.25
$ws.Range("L122").Value = 7085.9997
$ws.Range("M122").Value = -1989.25
$ws.Range("N122").Value = -11985.9997

# Row 132
$ws.Range("H132").Value = 1817.8
$ws.Range("I132").Value = 1500
$ws.Range("J132").Value = 1897.25
$ws.Range("K132").Value = 13500
$ws.Range("L132").Value = 17075.25
$ws.Range("M132").Value = -10970
$ws.Range("N132").Value = -22135.25

$ws = $wb.Worksheets.Item("GSM")
# Row 70
$ws.Range("H70").Value = 111114450
$ws.Range("I70").Value = 166669170
$ws.Range("J70").Value = 4999
$ws.Range("K70").Value = 166669170
$ws.Range("L70").Value = 4999
$ws.Range("M70").Value = -166668900
$ws.Range("N70").Value = -5539

# Row 73
$ws.Range("H73").Value = 111114450
$ws.Range("I73").Value = 166669170
$ws.Range("J73").Value = 4999
$ws.Range("K73").Value = 166669170
$ws.Range("L73").Value = 4999
$ws.Range("M73").Value = -166668234
$ws.Range("N73").Value = -6871

# Row 80
$ws.Range("H80").Value = 2004.5
$ws.Range("I80").Value = 2004.5
$ws.Range("J80").Value = 0
$ws.Range("K80").Value = 2004.5
$ws.Range("L80").Value = 0
$ws.Range("M80").Value = -1006.5
$ws.Range("N80").ClearContents()

# Row 83
$ws.Range("H83").Value = 2004.5
$ws.Range("I83").Value = 2004.5
$ws.Range("J83").Value = 0
$ws.Range("K83").Value = 10022.5
$ws.Range("L83").Value = 0
$ws.Range("M83").Value = -5030.5
$ws.Range("N83").ClearContents()

# Row 97
$ws.Range("H97").Value = 1577.3334
$ws.Range("I97").Value = 949.5
$ws.Range("K97").Value = 949.5
$ws.Range("M97").Value = -453.5

# Row 113
$ws.Range("H113").Value = 4997.5
$ws.Range("I113").Value = 4997.5
$ws.Range("K113").Value = 4997.5
$ws.Range("M113").Value = -2827.5

# Row 132
$ws.Range("H132").Value = 2332
$ws.Range("I132").Value = 2998
$ws.Range("J132").Value = 1999
$ws.Range("K132").Value = 8994
$ws.Range("L132").Value = 5997
$ws.Range("M132").Value = -6464
$ws.Range("N132").Value = -11057

$ws = $wb.Worksheets.Item("LTW")
# Row 7
$ws.Range("H7").Value = 21582.75
$ws.Range("I7").Value = 21582.75
$ws.Range("K7").Value = 21582.75
$ws.Range("M7").Value = -21470.75

# Row 40
$ws.Range("H40").Value = 9999.5
$ws.Range("J40").Value = 9999
$ws.Range("L40").Value = 9999
$ws.Range("N40").Value = -10271

# Row 93
$ws.Range("H93").Value = 3000
$ws.Range("I93").Value = 3000
$ws.Range("K93").Value = 3000
$ws.Range("M93").Value = -1752

# Row 126
$ws.Range("H126").Value = 21582.75
$ws.Range("I126").Value = 21582.75
$ws.Range("K126").Value = 64748.25
$ws.Range("M126").Value = -62278.25

$ws = $wb.Worksheets.Item("WVR")
# Row 81
$ws.Range("H81").Value = 1196.3334
$ws.Range("I81").Value = 1196.3334
$ws.Range("K81").Value = 2392.6668
$ws.Range("M81").Value = -1331.6668

# Row 84
$ws.Range("H84").Value = 1196.3334
$ws.Range("I84").Value = 1196.3334
$ws.Range("K84").Value = 11963.334
$ws.Range("M84").Value = -6659.333999999999

# Row 96
$ws.Range("H96").Value = 2596.75
$ws.Range("I96").Value = 2000
$ws.Range("J96").Value = 2795.6667
$ws.Range("K96").Value = 2000
$ws.Range("L96").Value = 2795.6667
$ws.Range("M96").Value = -627
$ws.Range("N96").Value = -5541.6667

# Row 100
$ws.Range("H100").Value = 1487.5
$ws.Range("J100").Value = 1475
$ws.Range("L100").Value = 2950
$ws.Range("N100").Value = -4032

# Row 136
$ws.Range("H136").Value = 6611.1333
$ws.Range("I136").Value = 4842.7144
$ws.Range("K136").Value = 14528.1432
$ws.Range("M136").Value = -11978.1432

